# Edit script applying the diff described in the commit:
# "added concepts for 2 and 3. 1 finished (without citations)"

$d = $word.ActiveDocument

function ReplaceText($findText, $replaceText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($found) {
        $rng.Text = $replaceText
    }
    return $found
}

# --- Paragraph 2: intro paragraph ---------------------------------------
ReplaceText "Системите за контрол на достъп са неделима част от всяка сграда" `
            "Системите за контрол са неделима част от инфраструктурата на почти всяка бизнес сграда"

ReplaceText "‚интелигентни‘" "IP базирани"

# --- Paragraph 3: "Целта..." ---------------------------------------------
ReplaceText "‚интелигентна‘ end-to-end система" "IP система"

ReplaceText " базиран четец и сървърна част за администриране, и осъществяване" `
            " базиран четец, сървърна част за администриране и осъществяване"

ReplaceText "идентификационния процес" "идентификационен процес"

ReplaceText "чрез Wi-fi комуникация." "чрез Wi-Fi комуникация между четеца и сървъра."

# --- Paragraph 5: Arduino board bullet -----------------------------------
ReplaceText "Ардуино борд (англ: Arduino board)" "Ардуино борд (англ.: Arduino board)"

# --- Paragraph 17: Wi-fi communication bullet -----------------------------
ReplaceText "Модул за Wi-fi комуникация" "Модул за Wi-Fi комуникация"

Write-Output "done"
